$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 377, pushing existing rows 377:392 down to 378:393.
$ws.Rows.Item(377).Insert()

# Populate the newly inserted row 377 with the new weekly price record.
$ws.Range("A377").Value = 4
$ws.Range("B377").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C377").Value = "Los Lagos"
$ws.Range("D377").Value = 45008
$ws.Range("E377").Value = 10
$ws.Range("F377").Value = 100112021
$ws.Range("G377").Value = "Ají"
$ws.Range("H377").Value = "Inferno"
$ws.Range("I377").Value = "Primera"
$ws.Range("J377").Value = 60
$ws.Range("K377").Value = 29000
$ws.Range("L377").Value = 29000
$ws.Range("M377").Value = 29000
$ws.Range("N377").Value = "$/caja 15 kilos"
$ws.Range("O377").Value = "Provincia de Quillota"
$ws.Range("P377").Value = 1933
$ws.Range("Q377").Value = 15
$ws.Range("R377").Value = "Hortaliza"

# Preserve the date-format style used by the other rows in column D.
$ws.Range("D377").NumberFormat = $ws.Range("D378").NumberFormat
